# "Merge Redhat and Ubuntu smoketesting"
#
# The "OS instructions" sheet previously had separate "Redhat" and "Ubuntu"
# rows (rows 4 and 5) that both pointed readers to the same tar.xz / conda
# install instructions. This merges them into a single "Linux" row (row 4),
# reusing the Linux tar.xz/conda instructions as the body text, and clears
# out the now-unused row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OS instructions")
$ws.Activate()

# Row 4: Redhat -> Linux, with the tar.xz/conda instructions (previously in B5)
$ws.Range("A4").Value = "Linux"
$linuxInstructions = "* To install the tar.xz package for Linux, run ``(sudo) tar -xJf mantid-VA.B.C.tar.xz`` in a terminal and it will unzip the package in your current working directory. `n* To install via conda:`n  - Use Intel Conda and make sure conda-forge is added to channels`n  - In terminal, create a new empty environment and activate it`n  - run ``conda install -c ""mantid/label/vA.B.C-rc1"" mantidworkbench`` , where A.B.C is the release version.`n"
$ws.Range("B4").Value = $linuxInstructions
$ws.Rows.Item(4).RowHeight = 132

# Row 5: clear out the old Ubuntu row, now unused
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Rows.Item(5).RowHeight = 20

# Update selection to reflect where the editor left off
$ws.Range("B7").Select() | Out-Null
